# Fixed glitch in data sample selection code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the mis-typed "other_crab" description
$ws.Range("C7").Value = "crabs, various sp"

# Restore the active cell/selection left by the editor
$ws.Activate()
$ws.Range("G13").Select()
